$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.037.82'
$ws.Range('E2').Value = '  +1.00%  '
$ws.Range('D3').Value = '2.274.84'
$ws.Range('E3').Value = '  +1.99%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '302.05'
$ws.Range('E5').Value = '  -1.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range('D6').Value = '94.72'
$ws.Range('E6').Value = '  +0.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range('D7').Value = '0.566'
$ws.Range('E7').Value = '  -1.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range('D8').Value = '0.997'
$ws.Range('E8').Value = '  -0.86%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range('D9').Value = '0.512'
$ws.Range('E9').Value = '  -0.88%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range('D10').Value = '34.29'
$ws.Range('E10').Value = '  -1.29%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range('D11').Value = '0.0793'
$ws.Range('E11').Value = '  -0.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range('D12').Value = '7.23'
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('E13').Value = '  -0.90%  '
$ws.Range('D14').Value = '2.613.18'
$ws.Range('E14').Value = '  +1.41%  '
$ws.Range('D15').Value = '2.265.42'
$ws.Range('E15').Value = '  +1.83%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range('D16').Value = '13.66'
$ws.Range('E16').Value = '  +1.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range('D17').Value = '0.803'
$ws.Range('E17').Value = '  -3.63%  '
$ws.Range('D18').Value = '44.892.05'
$ws.Range('E18').Value = '  +1.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range('D19').Value = '13.09'
$ws.Range('E19').Value = '  +9.81%  '
$ws.Range('E20').Value = '  -1.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range('D21').Value = '6.07'
$ws.Range('E21').Value = '  -2.57%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range('D22').Value = '65.71'
$ws.Range('E22').Value = '  +0.73%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range('D23').Value = '239.27'
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('E24').Value = '  -1.74%  '
$ws.Range('E25').Value = '  -0.30%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range('D26').Value = '1.90'
$ws.Range('E26').Value = '  -3.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range('D27').Value = '41.43'
$ws.Range('E27').Value = '  +10.21%  '
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range('D29').Value = '9.63'
$ws.Range('E29').Value = '  -1.46%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range('D30').Value = '19.68'
$ws.Range('E30').Value = '  -1.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range('D31').Value = '152.54'
$ws.Range('E31').Value = '  +1.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range('D32').Value = '5.57'
$ws.Range('E32').Value = '  -6.66%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range('D33').Value = '0.0794'
$ws.Range('E33').Value = '  +0.51%  '
$ws.Range('E34').Value = '  -2.28%  '
$ws.Range('E35').Value = '  -2.74%  '
$ws.Range('E36').Value = '  -1.09%  '
$ws.Range('E37').Value = '  -3.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range('D38').Value = '1.78'
$ws.Range('E38').Value = '  -2.76%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range('D39').Value = '3.97'
$ws.Range('E39').Value = '  +6.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range('D40').Value = '0.0311'
$ws.Range('E40').Value = '  +3.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range('D41').Value = '3.24'
$ws.Range('E41').Value = '  -3.53%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range('D42').Value = '13.72'
$ws.Range('E42').Value = '  -8.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  -0.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range('D44').Value = '1.94'
$ws.Range('E44').Value = '  +12.79%  '
$ws.Range('D45').Value = '1.749.09'
$ws.Range('E45').Value = '  -4.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range('D46').Value = '0.195'
$ws.Range('E46').Value = '  +3.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range('D47').Value = '76.75'
$ws.Range('E47').Value = '  -4.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range('D48').Value = '69.53'
$ws.Range('E48').Value = '  +1.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range('D49').Value = '95.64'
$ws.Range('E49').Value = '  -2.67%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '53.78'
$ws.Range('E50').Value = '  -0.46%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range('D51').Value = '4.70'
$ws.Range('E51').Value = '  -3.28%  '
